$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    $ws.Range($addr).Formula = "'" + $val
    $ws.Range($addr).Style = "Normal"
}

Set-TextCell "D2" "37.451.16"
Set-TextCell "E2" "  +1.22%  "
Set-TextCell "D3" "2.031.82"
Set-TextCell "E3" "  +0.46%  "
Set-TextCell "E5" "  +1.56%  "
Set-TextCell "D6" "0.612"
Set-TextCell "E6" "  +0.80%  "
Set-TextCell "E7" "  +0.00%  "
Set-TextCell "D8" "55.95"
Set-TextCell "E8" "  +1.73%  "
Set-TextCell "D9" "0.381"
Set-TextCell "E9" "  +0.18%  "
Set-TextCell "D10" "0.0799"
Set-TextCell "E10" "  +1.37%  "
Set-TextCell "E11" "  -1.24%  "
Set-TextCell "D12" "2.334.14"
Set-TextCell "E12" "  +0.37%  "
Set-TextCell "D13" "14.37"
Set-TextCell "E13" "  +0.46%  "
Set-TextCell "D14" "20.29"
Set-TextCell "E14" "  -0.49%  "
Set-TextCell "D15" "0.741"
Set-TextCell "E15" "  -0.38%  "
Set-TextCell "D16" "5.21"
Set-TextCell "E16" "  +1.13%  "
Set-TextCell "D17" "2.036.02"
Set-TextCell "E17" "  +0.14%  "
Set-TextCell "D18" "37.409.79"
Set-TextCell "E18" "  +1.18%  "
Set-TextCell "D19" "6.19"
Set-TextCell "E19" "  -0.92%  "
Set-TextCell "D20" "68.99"
Set-TextCell "E20" "  +0.32%  "
Set-TextCell "D21" "0.0₃0824"
Set-TextCell "E21" "  +0.00%  "
Set-TextCell "D22" "223.53"
Set-TextCell "E22" "  -1.08%  "
Set-TextCell "E23" "  -0.02%  "
Set-TextCell "E24" "  +1.59%  "
Set-TextCell "E25" "  +2.77%  "
Set-TextCell "D26" "165.04"
Set-TextCell "E26" "  -0.33%  "
Set-TextCell "D27" "9.12"
Set-TextCell "E27" "  -1.33%  "
Set-TextCell "E28" "  +5.63%  "
Set-TextCell "D29" "18.74"
Set-TextCell "E29" "  +0.13%  "
Set-TextCell "E30" "  -1.18%  "
Set-TextCell "E31" "  +0.69%  "
Set-TextCell "D32" "4.48"
Set-TextCell "E32" "  +0.23%  "
Set-TextCell "D33" "0.0606"
Set-TextCell "E33" "  -1.84%  "
Set-TextCell "E34" "  +0.90%  "
Set-TextCell "E35" "  +8.65%  "
Set-TextCell "D36" "2.31"
Set-TextCell "E36" "  -1.68%  "
Set-TextCell "B37" "RenderToken"
Set-TextCell "C37" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextCell "D37" "3.25"
Set-TextCell "E37" "  +3.00%  "
Set-TextCell "B38" "THORChain"
Set-TextCell "C38" "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
Set-TextCell "D38" "5.76"
Set-TextCell "E38" "  +8.63%  "
Set-TextCell "E39" "  +0.14%  "
Set-TextCell "D40" "1.474.99"
Set-TextCell "E40" "  -0.92%  "
Set-TextCell "D41" "0.0214"
Set-TextCell "E41" "  -1.52%  "
Set-TextCell "D42" "0.0930"
Set-TextCell "E42" "  +0.46%  "
Set-TextCell "E43" "  +1.43%  "
Set-TextCell "D44" "94.93"
Set-TextCell "E44" "  -0.42%  "
Set-TextCell "D45" "16.36"
Set-TextCell "E45" "  -5.13%  "
Set-TextCell "D46" "4.19"
Set-TextCell "E46" "  +15.94%  "
Set-TextCell "E47" "  -2.13%  "
Set-TextCell "E48" "  +0.19%  "
Set-TextCell "B49" "MXToken"
Set-TextCell "C49" "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextCell "D49" "2.94"
Set-TextCell "E49" "  +0.67%  "
Set-TextCell "B50" "FraxShare"
Set-TextCell "C50" "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextCell "D50" "7.09"
Set-TextCell "E50" "  -3.82%  "
Set-TextCell "D51" "2.221.72"
Set-TextCell "E51" "  +0.37%  "
